$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet (sheet4.xml): insert a new blank column
# before column N. This shifts the old N:P columns ("Late" label +
# values, and the "Outstanding" label + values) one column to the
# right, becoming O:Q, and leaves a new blank column N behind.
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Range("N1").EntireColumn.Insert()

# "Repayment Schedule" becomes the selected/active sheet (activeTab
# moves from the "Transactions" sheet to this one), with S5 selected.
$ws.Activate()
$ws.Range("S5").Select()
